$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.614.87'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.564.06'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.22'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.509'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.93'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +5.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.245'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0587'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.787.21'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.561.25'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.669.26'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.516'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.51'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '227.78'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0680'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.43%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.64'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.05%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.77'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0458'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.93%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.18'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.399.83'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.02%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.69%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.79%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.30'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.66'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.21%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.03%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.766'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.84'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.700.22'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.868'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '84.76'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.60%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.86%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.44%  '
